# Apply the "10Th - MB for single stock and added new group" update.
#
# Summary of the change:
#  - Insert 3 new weekly snapshot columns (B:D) in front of the existing
#    B:E block, pushing the old B/C/D/E columns to E/F/G/H.
#  - New header cells: B1=Jun_27, C1=Jun_26, D1=Jun_26 (E1:H1 keep the old
#    Jun_17 / Jun_15 / Jun_13 / Jun_10 headers, now shifted right).
#  - The new B:D columns are filled with "UN" (unchanged) for every
#    existing analyst row (2-27), matching how E:G already look before the
#    insert.
#  - Two new analyst rows are appended: Benchmark (28) and Evercore ISI (29),
#    each with "UN" in B:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert three new columns before column B -------------------------
# This shifts the existing B:E columns (and their widths) to E:H.
$ws.Range("B1:D1").EntireColumn.Insert()

# --- 2. New header row values for the freshly inserted columns -----------
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- 3. Append the two new analyst rows -----------------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

$ws.Range("B1").Value = "Jun_27"

# --- 4. Fill the new columns with "UN" for every existing data row -------
$lastRow = 27
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("B" + $r).Value = "UN"
    $ws.Range("C" + $r).Value = "UN"
    $ws.Range("D" + $r).Value = "UN"
}
